$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 25 de Marzo de 2020 a las 15:46'

$ws.Cells.Item(8, 2).Value = 35704
$ws.Cells.Item(8, 3).Value = 2713
$ws.Cells.Item(8, 5).Value = 31983

$ws.Cells.Item(19, 4).Value = 185
$ws.Cells.Item(19, 5).Value = 2581

$ws.Cells.Item(20, 5).Value = 2466
$ws.Cells.Item(20, 7).Value = 4
$ws.Cells.Item(20, 8).Value = 44

$ws.Cells.Item(30, 5).Value = 1117
$ws.Cells.Item(30, 7).Value = 1
$ws.Cells.Item(30, 8).Value = 3

$ws.Cells.Item(46, 1).Value = 'Singapur'
$ws.Cells.Item(46, 2).Value = 631
$ws.Cells.Item(46, 3).Value = 73
$ws.Cells.Item(46, 4).Value = 156
$ws.Cells.Item(46, 5).Value = 473
$ws.Cells.Item(46, 6).Value = 17
$ws.Cells.Item(46, 8).Value = 2

$ws.Cells.Item(47, 1).Value = 'India'
$ws.Cells.Item(47, 2).Value = 606
$ws.Cells.Item(47, 3).Value = 70
$ws.Cells.Item(47, 4).Value = 42
$ws.Cells.Item(47, 5).Value = 554
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 8).Value = 10

$ws.Cells.Item(58, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(58, 2).Value = 392
$ws.Cells.Item(58, 3).Value = 80
$ws.Cells.Item(58, 4).Value = 3
$ws.Cells.Item(58, 5).Value = 379
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 10

$ws.Cells.Item(59, 1).Value = 'Argentina'
$ws.Cells.Item(59, 2).Value = 387
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 52
$ws.Cells.Item(59, 5).Value = 328
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 8).Value = 7

$ws.Cells.Item(60, 1).Value = 'Serbia'
$ws.Cells.Item(60, 2).Value = 384
$ws.Cells.Item(60, 3).Value = 81
$ws.Cells.Item(60, 4).Value = 15
$ws.Cells.Item(60, 5).Value = 365
$ws.Cells.Item(60, 6).Value = 21
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 4

$ws.Cells.Item(61, 1).Value = 'Colombia'
$ws.Cells.Item(61, 2).Value = 378
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 6
$ws.Cells.Item(61, 5).Value = 369
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 3

$ws.Cells.Item(62, 1).Value = 'Irak'
$ws.Cells.Item(62, 2).Value = 346
$ws.Cells.Item(62, 3).Value = 30
$ws.Cells.Item(62, 4).Value = 103
$ws.Cells.Item(62, 5).Value = 214
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = 29

$ws.Cells.Item(63, 1).Value = 'Libano'
$ws.Cells.Item(63, 3).Value = 15
$ws.Cells.Item(63, 4).Value = 8
$ws.Cells.Item(63, 5).Value = 321
$ws.Cells.Item(63, 6).Value = 4
$ws.Cells.Item(63, 8).Value = 4

$ws.Cells.Item(64, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(64, 2).Value = 333
$ws.Cells.Item(64, 3).Value = 85
$ws.Cells.Item(64, 4).Value = 52
$ws.Cells.Item(64, 5).Value = 279
$ws.Cells.Item(64, 6).Value = 2
$ws.Cells.Item(64, 8).Value = 2

$ws.Cells.Item(80, 1).Value = 'Tunez'
$ws.Cells.Item(80, 2).Value = 171
$ws.Cells.Item(80, 3).Value = 57
$ws.Cells.Item(80, 4).Value = 2
$ws.Cells.Item(80, 5).Value = 165
$ws.Cells.Item(80, 6).Value = 11
$ws.Cells.Item(80, 8).Value = 4

$ws.Cells.Item(81, 1).Value = 'Marruecos'
$ws.Cells.Item(81, 2).Value = 170
$ws.Cells.Item(81, 4).Value = 6
$ws.Cells.Item(81, 5).Value = 159
$ws.Cells.Item(81, 8).Value = 5

$ws.Cells.Item(82, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(82, 2).Value = 168
$ws.Cells.Item(82, 4).Value = 2
$ws.Cells.Item(82, 5).Value = 163
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 8).Value = 3

$ws.Cells.Item(83, 1).Value = 'Jordania'
$ws.Cells.Item(83, 2).Value = 154
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 1
$ws.Cells.Item(83, 5).Value = 153
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 8).Value = 0

$ws.Cells.Item(84, 1).Value = 'Albania'
$ws.Cells.Item(84, 2).Value = 146
$ws.Cells.Item(84, 3).Value = 23
$ws.Cells.Item(84, 8).Value = 5

$ws.Cells.Item(85, 1).Value = 'Vietnam'
$ws.Cells.Item(85, 2).Value = 141
$ws.Cells.Item(85, 3).Value = 7
$ws.Cells.Item(85, 4).Value = 17
$ws.Cells.Item(85, 5).Value = 124
$ws.Cells.Item(85, 6).Value = 3

$ws.Cells.Item(86, 1).Value = 'Islas Feroe'
$ws.Cells.Item(86, 2).Value = 132
$ws.Cells.Item(86, 3).Value = 10
$ws.Cells.Item(86, 4).Value = 38
$ws.Cells.Item(86, 5).Value = 94
$ws.Cells.Item(86, 6).Value = 2

$ws.Cells.Item(87, 1).Value = 'Malta'
$ws.Cells.Item(87, 2).Value = 129
$ws.Cells.Item(87, 3).Value = 19
$ws.Cells.Item(87, 5).Value = 127
$ws.Cells.Item(87, 6).Value = 1
$ws.Cells.Item(87, 8).Value = 0

$ws.Cells.Item(88, 1).Value = 'Moldavia'
$ws.Cells.Item(88, 2).Value = 125
$ws.Cells.Item(88, 4).Value = 2
$ws.Cells.Item(88, 5).Value = 122
$ws.Cells.Item(88, 6).Value = 20
$ws.Cells.Item(88, 8).Value = 1

$ws.Cells.Item(89, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(89, 2).Value = 124
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 3
$ws.Cells.Item(89, 5).Value = 118
$ws.Cells.Item(89, 6).Value = 3
$ws.Cells.Item(89, 8).Value = 3

$ws.Cells.Item(90, 1).Value = 'Ucrania'
$ws.Cells.Item(90, 2).Value = 116
$ws.Cells.Item(90, 3).Value = 14
$ws.Cells.Item(90, 4).Value = 1
$ws.Cells.Item(90, 5).Value = 111
$ws.Cells.Item(90, 7).Value = 1

$ws.Cells.Item(91, 1).Value = 'Burkina Faso'
$ws.Cells.Item(91, 2).Value = 114
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 7
$ws.Cells.Item(91, 5).Value = 103
$ws.Cells.Item(91, 7).Value = 0

$ws.Cells.Item(99, 5).Value = 75
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 2

$ws.Cells.Item(114, 1).Value = 'Mauricio'
$ws.Cells.Item(114, 3).Value = 6
$ws.Cells.Item(114, 6).Value = 1

$ws.Cells.Item(115, 1).Value = 'Cuba'
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 1
$ws.Cells.Item(115, 6).Value = 2
$ws.Cells.Item(115, 8).Value = 1

$ws.Cells.Item(116, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(116, 3).Value = 3
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 8).Value = 2

$ws.Cells.Item(123, 1).Value = 'Mayotte'
$ws.Cells.Item(123, 3).Value = 0

$ws.Cells.Item(124, 1).Value = 'Honduras'
$ws.Cells.Item(124, 3).Value = 6

$ws.Cells.Item(128, 1).Value = 'Polinesia Francesa'

$ws.Cells.Item(129, 1).Value = 'Kenia'

$ws.Cells.Item(141, 1).Value = 'Uganda'
$ws.Cells.Item(141, 3).Value = 5

$ws.Cells.Item(142, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(142, 3).Value = 4

$ws.Cells.Item(144, 1).Value = 'Tanzania'

$ws.Cells.Item(145, 1).Value = 'Etiopia'
$ws.Cells.Item(145, 3).Value = 0

$ws.Cells.Item(146, 1).Value = 'Zambia'
$ws.Cells.Item(146, 3).Value = 9

$ws.Cells.Item(152, 1).Value = 'Seychelles'

$ws.Cells.Item(154, 1).Value = 'Haiti'

$ws.Cells.Item(155, 1).Value = 'Dominica'

$ws.Cells.Item(161, 1).Value = 'Islas Caimanes'

$ws.Cells.Item(162, 1).Value = 'Curazao'

$ws.Cells.Item(164, 1).Value = 'Bahamas'
$ws.Cells.Item(164, 4).Value = 1
$ws.Cells.Item(164, 8).Value = 0

$ws.Cells.Item(165, 1).Value = 'Guyana'
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 8).Value = 1

$ws.Cells.Item(167, 1).Value = 'Congo'

$ws.Cells.Item(170, 1).Value = 'Suazilandia'

$ws.Cells.Item(171, 1).Value = 'Cabo Verde'
$ws.Cells.Item(171, 3).Value = 1
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 8).Value = 1

$ws.Cells.Item(172, 1).Value = 'Birmania'
$ws.Cells.Item(172, 2).Value = 3
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 1).Value = 'Mozambique'

$ws.Cells.Item(175, 1).Value = 'Liberia'

$ws.Cells.Item(176, 1).Value = 'San Bartolome'

$ws.Cells.Item(178, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(178, 3).Value = 0

$ws.Cells.Item(179, 1).Value = 'Angola'

$ws.Cells.Item(181, 1).Value = 'Laos'
$ws.Cells.Item(181, 3).Value = 1

$ws.Cells.Item(182, 1).Value = 'Nepal'
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 5).Value = 2

$ws.Cells.Item(184, 1).Value = 'Sudan'

$ws.Cells.Item(185, 1).Value = 'Zimbabue'

$ws.Cells.Item(186, 1).Value = 'Nicaragua'

$ws.Cells.Item(187, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(187, 3).Value = 2

$ws.Cells.Item(189, 1).Value = 'Mauritania'

$ws.Cells.Item(190, 1).Value = 'Butan'
$ws.Cells.Item(190, 3).Value = 0

$ws.Cells.Item(191, 1).Value = 'Mali'

$ws.Cells.Item(192, 1).Value = 'Siria'

$ws.Cells.Item(193, 1).Value = 'Eritrea'

$ws.Cells.Item(194, 1).Value = 'Timor Oriental'

$ws.Cells.Item(196, 1).Value = 'Belice'

$ws.Cells.Item(197, 1).Value = 'Montserrat'

$ws.Cells.Item(198, 1).Value = 'Papua Nueva Guinea'

$ws.Cells.Item(200, 1).Value = 'Somalia'

$ws.Cells.Item(201, 1).Value = 'San Vicente y las Granadinas'

$ws.Cells.Item(202, 1).Value = 'Granada'
